$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text instead of coercing to a numeric cell (matches
# the original inline-string / text semantics of these cells).
$textCells = @("D5", "D10", "D14", "D15", "D18", "D20", "D26", "D29", "D30", "D31", "D33", "D34", "D35", "D37", "D41", "D42", "D46", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '26.015.44'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.633.82'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '214.73'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('D10').Value = '18.38'
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').Value = '1.860.80'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '1.630.33'
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').Value = '4.17'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = '0.525'
$ws.Range('E15').Value = '  -3.55%  '
$ws.Range('D16').Value = '25.994.31'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = '0.0₃0742'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').Value = '61.42'
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '190.61'
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('E21').Value = '  -2.40%  '
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('D26').Value = '143.56'
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -2.12%  '
$ws.Range('D29').Value = '15.17'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('D31').Value = '0.0481'
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('D33').Value = '3.13'
$ws.Range('E33').Value = '  -4.94%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.40'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('D36').Value = '1.132.23'
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '0.860'
$ws.Range('E37').Value = '  -5.18%  '
$ws.Range('E39').Value = '  -4.61%  '
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('D41').Value = '98.33'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('D42').Value = '0.774'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('D43').Value = '1.771.00'
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('E44').Value = '  -5.20%  '
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('D46').Value = '54.76'
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '1.49'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').Value = '7.48'
$ws.Range('E51').Value = '  -3.40%  '
